# Add Slack IDs for leaders ("lider_slack_id" / "slack_lider" columns)
# to the mapeamento_final and resumo_lideres sheets.

$wb = $excel.ActiveWorkbook

# Map of normalized leader name -> Slack user id.
$leaderSlackIds = @{
    "Romulo Jose Santos Lisboa"             = "U07LSKN7SNL"
    "Alberto Luiz Marinho Batista"          = "U07KXEJU338"
    "Joao Antonio Tavares Santos"           = "U07LP4JSN9K"
    "Jonathan Henrique da Conceição Silva"  = "U07L4D3EWJW"
    "Carlos Eduardo Silva De Oliveira"      = "U0895CZ8HU7"
    "Leidiane Souza"                        = "U07KX76F7D4"
    "Erick Café Santos Júnior"              = "U07KPE840MD"
    "Ana Clara de Matos Chagas"             = "U08F9KK0AAG"
    "Kemilly Rafaelly Souza Silva"          = "U087HDEARA9"
    "Maria Taciane Pereira Barbosa"         = "U07L6EAUS75"
    "Mariane Santos Sousa"                  = "U088B372R40"
    "Michaell Jean Nunes De Carvalho"       = "U07P692F1FB"
    "Rafaela Alves Mendes"                  = "U07KP9J5BLP"
    "Suzana Martins Tavares"                = "U09F9LWM6MC"
    "Ravy Thiago Vieira Da Silva"           = "U07Q8NT7J1Y"
}

# "Leidiane Souza.1" is a deduplicated-name variant of "Leidiane Souza" that
# appears further down the mapeamento_final sheet; it shares her Slack id.
$leaderSlackIdsMap = $leaderSlackIds.Clone()
$leaderSlackIdsMap["Leidiane Souza.1"] = "U07KX76F7D4"

# --- Sheet "mapeamento_final": column F (6) = lider_slack_id, rows 2-78 ---
$wsMap = $wb.Worksheets.Item("mapeamento_final")

for ($row = 2; $row -le 78; $row++) {
    $leaderName = $wsMap.Cells.Item($row, 1).Value2
    if ($leaderSlackIdsMap.ContainsKey($leaderName)) {
        $wsMap.Cells.Item($row, 6).Value = $leaderSlackIdsMap[$leaderName]
    }
}

# --- Sheet "resumo_lideres": column C (3) = slack_lider, rows 2-17 ---
$wsResumo = $wb.Worksheets.Item("resumo_lideres")

for ($row = 2; $row -le 17; $row++) {
    $leaderName = $wsResumo.Cells.Item($row, 1).Value2
    if ($leaderSlackIds.ContainsKey($leaderName)) {
        $wsResumo.Cells.Item($row, 3).Value = $leaderSlackIds[$leaderName]
    }
}
